# Commit: Changed Roxie interpolation method and plot
# Update the recomputed "Roxie" interpolation result columns (R:AF)
# for data rows 2-6 on Sheet1 to the new interpolation-method values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")


# Row 2
$ws.Range("R2").Value = -0.2988425617549172
$ws.Range("S2").Value = -9.727776421031244
$ws.Range("T2").Value = -23.6164547842054
$ws.Range("U2").Value = 0.07297155402910539
$ws.Range("V2").Value = 3.163732635147427
$ws.Range("W2").Value = 0.03455639280639029
$ws.Range("X2").Value = -3.749163636960444
$ws.Range("Y2").Value = -0.03719803425637391
$ws.Range("Z2").Value = -1.430227016780397
$ws.Range("AA2").Value = 0.004339350245012328
$ws.Range("AB2").Value = -1.813934712341678
$ws.Range("AC2").Value = 0.03295322185671454
$ws.Range("AD2").Value = 7.736624931039496
$ws.Range("AE2").Value = -0.02138943740816093
$ws.Range("AF2").Value = -4.668774128269343

# Row 3
$ws.Range("R3").Value = 3883.249378802444
$ws.Range("S3").Value = -4.310146059960727
$ws.Range("T3").Value = -102.2295270416021
$ws.Range("U3").Value = 0.2098808488041859
$ws.Range("V3").Value = -41.24458552371551
$ws.Range("W3").Value = 0.009078976583652964
$ws.Range("X3").Value = -4.916793741932531
$ws.Range("Y3").Value = 0.006673221595687773
$ws.Range("Z3").Value = -1.727976748465495
$ws.Range("AA3").Value = 0.002233848600394632
$ws.Range("AB3").Value = -0.5432594962308495
$ws.Range("AC3").Value = -0.004737495432041968
$ws.Range("AD3").Value = -0.6079292103552552
$ws.Range("AE3").Value = -0.004324656446714305
$ws.Range("AF3").Value = 0.5943259143989457

# Row 4
$ws.Range("R4").Value = 9999.54449330356
$ws.Range("S4").Value = -0.00141798797441333
$ws.Range("T4").Value = -13.01506571148024
$ws.Range("U4").Value = -0.0003021302918249709
$ws.Range("V4").Value = 0.2767574585328175
$ws.Range("W4").Value = -0.00005063763162227398
$ws.Range("X4").Value = 1.78449158866461
$ws.Range("Y4").Value = 0.000002816113789196996
$ws.Range("Z4").Value = 2.576922215211461
$ws.Range("AA4").Value = -0.000003348889359806453
$ws.Range("AB4").Value = 3.086940170692143
$ws.Range("AC4").Value = 0.000001230594817298837
$ws.Range("AD4").Value = -1.695117678733872
$ws.Range("AE4").Value = -0.000002406476643418612
$ws.Range("AF4").Value = 0.3234253558622677

# Row 5
$ws.Range("R5").Value = 3924.583743737668
$ws.Range("S5").Value = -0.2664069153075049
$ws.Range("T5").Value = 732.0905336348602
$ws.Range("U5").Value = -0.01983842572494893
$ws.Range("V5").Value = 357.9478538335806
$ws.Range("W5").Value = 0.001331799116434746
$ws.Range("X5").Value = 124.1039117647692
$ws.Range("Y5").Value = -0.000581540885001387
$ws.Range("Z5").Value = 45.68089166192325
$ws.Range("AA5").Value = -0.00002131295993674912
$ws.Range("AB5").Value = 17.56577122913476
$ws.Range("AC5").Value = 0.0006067413513743717
$ws.Range("AD5").Value = 5.911530117808881
$ws.Range("AE5").Value = 0.0003945885488136774
$ws.Range("AF5").Value = 2.44918292339739

# Row 6
$ws.Range("R6").Value = 1.629919257347524
$ws.Range("S6").Value = -0.2241967543583217
$ws.Range("T6").Value = 13.36935793973882
$ws.Range("U6").Value = -0.01011433821218637
$ws.Range("V6").Value = -0.3463083933613519
$ws.Range("W6").Value = 0.009535958439886217
$ws.Range("X6").Value = 0.6722488571783672
$ws.Range("Y6").Value = 0.005900677588224554
$ws.Range("Z6").Value = 0.2557707902136752
$ws.Range("AA6").Value = 0.0006709541502853638
$ws.Range("AB6").Value = 0.4136810761555675
$ws.Range("AC6").Value = -0.008947630803538735
$ws.Range("AD6").Value = -1.197859966876656
$ws.Range("AE6").Value = -0.005015797054788405
$ws.Range("AF6").Value = 0.5390454768932538
